$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "Thyroid" entry as row 24 (sheet grows from A1:E23 to A1:E24).
$ws.Range("A24").Value = "Thyroid"
$ws.Range("B24").Value = "Isoechoic nodule with peripheral calcifications"
$ws.Range("C24").Value = "Clip 1 B-mode + Color"

# D24 gets the YouTube link as a hyperlink (same pattern used by every other
# row's "YouTube Link" column), with the built-in hyperlink cell style.
$ws.Hyperlinks.Add($ws.Range("D24"), "https://youtu.be/z_oaRVxRz5s ", "", "", "https://youtu.be/z_oaRVxRz5s ") | Out-Null
$ws.Range("D24").Style = "Collegamento ipertestuale"

# Restore the author's final selection noted in the saved workbook.
$ws.Range("D27").Select() | Out-Null
